$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 83, shifting existing rows 83:160 down to 84:161.
$ws.Rows("83:83").Insert()

# Populate the newly inserted row 83 with the new record.
$ws.Range("A83").Value = 10
$ws.Range("B83").Value = "Vega Modelo de Temuco"
$ws.Range("C83").Value = "La Araucanía"
$ws.Range("D83").Value = 44484
$ws.Range("E83").Value = 9
$ws.Range("F83").Value = 100112043
$ws.Range("G83").Value = "Pepino dulce"
$ws.Range("H83").Value = "Cultivar IV Región"
$ws.Range("I83").Value = "Segunda"
$ws.Range("J83").Value = 20
$ws.Range("K83").Value = 20000
$ws.Range("L83").Value = 20000
$ws.Range("M83").Value = 20000
$ws.Range("N83").Value = "$/bandeja 18 kilos"
$ws.Range("O83").Value = "Provincia de Limarí"
$ws.Range("P83").Value = 1111
$ws.Range("Q83").Value = 18
$ws.Range("R83").Value = "Hortaliza"
